# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# Both sheets contain identical data rows, and both need the same updates.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    2  = 11493
    3  = 10955
    5  = 6
    6  = 993
    7  = 116
    8  = 58
    10 = 37
    11 = 10618
    12 = 4104
    13 = 8
    16 = 36
    18 = 421
    19 = 11101
    20 = 10862
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
